$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(228, 1).Value = 'https://www.dakotanewsnow.com/2025/02/14/jackley-clarifies-focus-multi-state-section-504-lawsuit/'
$ws.Cells.Item(228, 2).Value = 'Jackley clarifies the focus of multi-state Section 504 lawsuit'

$ws.Cells.Item(229, 1).Value = 'https://www.dallasnews.com/news/2025/02/14/lawsuit-filed-friday-against-keller-isd-alleges-violations-of-voting-rights-act/'
$ws.Cells.Item(229, 2).Value = 'Lawsuit filed against Keller ISD alleges violations of Voting Rights Act'

$ws.Cells.Item(230, 1).Value = 'https://arkansasadvocate.com/briefs/20-red-states-including-arkansas-back-doge-in-lawsuit-challenging-access-to-treasury-system/'
$ws.Cells.Item(230, 2).Value = '20 red states, including Arkansas, back DOGE in lawsuit challenging access to Treasury system'

$ws.Cells.Item(231, 1).Value = 'https://lawnews.hofstra.edu/2025/02/14/prof-james-sample-explains-federal-lawsuit-against-new-yorks-green-light-law/'
$ws.Cells.Item(231, 2).Value = 'Prof. James Sample Explains Federal Lawsuit Against New York’s Green Light Law'

$ws.Cells.Item(232, 1).Value = 'https://thedaily.case.edu/laws-sharona-hoffman-discusses-an-insulin-lawsuit-recently-filed-by-the-city-of-columbus/'
$ws.Cells.Item(232, 2).Value = 'Law’s Sharona Hoffman discusses an insulin lawsuit recently filed by the City of Columbus'

$ws.Cells.Item(233, 1).Value = 'https://www.alexcityoutlook.com/elmore-county-lawyer-arrested-for-harassment/article_87d9737c-4b06-5dae-aa75-326147315dd6.html'
$ws.Cells.Item(233, 2).Value = 'Elmore County lawyer arrested for harassment'

$ws.Cells.Item(234, 1).Value = 'https://www.yahoo.com/entertainment/fla-lawyer-allegedly-smashed-plate-141907407.html'
$ws.Cells.Item(234, 2).Value = 'Fla. Lawyer Allegedly Smashed Plate Over Fellow Wedding Attendee’s Head When He Allowed Others to Cut in Buffet Line'

$ws.Cells.Item(235, 1).Value = 'https://www.whas11.com/video/news/local/indiana/attorney-general-threatening-legal-action-against-indianapolis-authorities-schools/417-bcc673c1-09e4-4abc-b8fc-738431a1ec29'
$ws.Cells.Item(235, 2).Value = 'Attorney general threatening legal action against Indianapolis authorities, schools'

$ws.Cells.Item(236, 1).Value = 'https://edmontonjournal.com/news/politics/outrageous-and-false-how-those-named-in-the-ahs-lawsuit-are-responding-to-the-allegations'
$ws.Cells.Item(236, 2).Formula = "=""'Outrageous and false': How those named in the AHS lawsuit are responding to the allegations"""

$ws.Cells.Item(237, 1).Value = 'https://www.billboard.com/music/rb-hip-hop/glorilla-bbl-glorious-tour-1235903399/'
$ws.Cells.Item(237, 2).Value = 'GloRilla Denies BBL Rumors as She Preps for ‘The Glorious Tour’'

$ws.Cells.Item(238, 1).Value = 'https://theprint.in/india/ed-seizes-rs-170-cr-worth-bank-deposits-in-probe-against-fraud-forex-trading-platform/2492261/'
$ws.Cells.Item(238, 2).Value = 'ED seizes Rs 170-cr worth bank deposits in probe against ‘fraud’ forex trading platform'

$ws.Cells.Item(239, 1).Value = 'https://www.notus.org/whitehouse/doge-posts-then-redacts-sensitive-hud-contract-data'
$ws.Cells.Item(239, 2).Value = 'DOGE Posts — Then Redacts — What Appears to Be Sensitive HUD Contract Data'

$ws.Cells.Item(240, 1).Value = 'https://www.fox44news.com/news/local-news/bell-county/tarver-elementary-teacher-placed-on-leave-amid-misconduct-allegations/'
$ws.Cells.Item(240, 2).Value = 'Tarver Elementary teacher placed on leave amid misconduct allegations - KWKT'

$ws.Cells.Item(241, 1).Value = 'https://mynbc15.com/news/local/mobile-mardi-gras-queen-accused-of-nearly-15m-embezzlement-scheme'
$ws.Cells.Item(241, 2).Value = 'Mobile Mardi Gras queen accused of nearly $1.5M embezzlement scheme'

$ws.Cells.Item(242, 1).Value = 'https://www.texomashomepage.com/news/crime/vernon-business-owner-accused-of-child-sex-crimes/'
$ws.Cells.Item(242, 2).Value = 'Vernon business owner accused of child sex crimes - KFDX'

$ws.Cells.Item(243, 1).Value = 'https://www.msn.com/en-us/news/crime/mass-daycare-co-owner-accused-of-assaulting-children-in-her-care/ar-AA1yyDPu'
$ws.Cells.Item(243, 2).Value = 'Mass. daycare co-owner accused of assaulting children in her care'

$ws.Cells.Item(244, 1).Value = 'https://www.m9.news/usa-news/b1-b2-visa-revoked-traveler-finds-out-at-airport/'
$ws.Cells.Item(244, 2).Value = 'B1/B2 Visa Revoked -Traveler Finds Out at Airport'

$ws.Cells.Item(245, 1).Value = 'https://www.bbc.com/news/articles/c9qj8gelgz8o'
$ws.Cells.Item(245, 2).Value = 'EHarley Street: Calls for inquiry into GP management ''scandal'''

$ws.Cells.Item(246, 1).Value = 'https://www.newsx.com/entertainment/netflix-addresses-controversy-surrounding-karla-sofia-gascon-amid-emilia-perez-scandal/'
$ws.Cells.Item(246, 2).Value = 'Netflix Addresses Controversy Surrounding Karla Sofía Gascón Amid ‘Emilia Pérez’ Scandal'

$ws.Range("B236").Copy() | Out-Null
$ws.Range("B236").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false
